$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card18")

# --- Row 13: previously-blank columns now logged as "nan" (text) ---
$nanCols13 = @("B","C","D","E","F","G","H","I","J","K","M","N")
foreach ($col in $nanCols13) {
    $ws.Range($col + "13").Value = "nan"
}
# O13 and P13 already contain their final values; L13 already has its date.

# --- Row 14: brand-new service event appended to the log ---
# A14 must be stored as TEXT "18" (not a number), matching the rest of column A.
$a14 = $ws.Range("A14")
$a14.NumberFormat = "@"
$a14.Value = "18"
$a14.ClearFormats()

# B14:K14 and M14:N14 stay blank (present but empty), like the template rows.
$blankCols14 = @("B","C","D","E","F","G","H","I","J","K","M","N")
foreach ($col in $blankCols14) {
    $cell = $ws.Range($col + "14")
    $cell.Value = ""
    $cell.ClearFormats()
}

$ws.Range("L14").Value = "13\8\2024"
$ws.Range("O14").Value = "تم سير مضرب dfk 25*1.5*974 flat"
$ws.Range("P14").Value = "تيم العمل"
